# Apply the edits described in the commit:
#  - Rename the three worksheets (translate Russian default names to English)
#  - Fix the date typo in the report header (04-10-2018 -> 05-10-2018)
#  - Move the remembered cell selection on every sheet from C15 to C26

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Rename worksheets
$ws1.Name = "Worksheet 1"
$ws2.Name = "Workshet 2"
$ws3.Name = "Workshet 3"

# Correct the date in the shared header text (cell A2, merged A2:D2) on every sheet
$newHeader = "Date: 05-10-2018 - Department: Sales department"
$ws1.Range("A2").Value = $newHeader
$ws2.Range("A2").Value = $newHeader
$ws3.Range("A2").Value = $newHeader

# Move the active-cell selection from C15 to C26 on every sheet,
# keeping the first worksheet as the active tab/selection at the end.
$ws1.Activate()
$ws1.Range("C26").Select()

$ws2.Activate()
$ws2.Range("C26").Select()

$ws3.Activate()
$ws3.Range("C26").Select()

$ws1.Activate()
